$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.097.09"
$ws.Range("E2").Value = "  -3.05%  "
$ws.Range("D3").Value = "3.325.68"
$ws.Range("E3").Value = "  -5.22%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "162.50"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -7.21%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("B8").Value = "LidoStakedEther"
$ws.Range("C8").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D8").Value = "3.321.06"
$ws.Range("E8").Value = "  -5.21%  "
$ws.Range("B9").Value = "XRP"
$ws.Range("C9").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.572"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -6.33%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.184"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.65"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -7.12%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.526"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -9.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "41.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -10.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000255"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.32%  "
$ws.Range("D15").Value = "3.861.89"
$ws.Range("E15").Value = "  -5.25%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "68.024.41"
$ws.Range("E16").Value = "  -3.31%  "
$ws.Range("B17").Value = "Polkadot"
$ws.Range("C17").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.64"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -7.78%  "
$ws.Range("D18").Value = "3.327.62"
$ws.Range("E18").Value = "  -5.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "553.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -9.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.117"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.01"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.796"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -9.03%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.34"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -7.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "89.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -7.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.34%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -6.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("B28").Value = "ImmutableX"
$ws.Range("C28").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.54%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "30.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -9.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -9.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -9.92%  "
$ws.Range("E32").Value = "  -7.24%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.64"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -11.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "579.84"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -9.59%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -10.16%  "
$ws.Range("E36").Value = "  +0.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "55.08"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("B38").Value = "Cosmos"
$ws.Range("C38").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -8.19%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0445"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.56%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0898"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -9.60%  "
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -18.95%  "
$ws.Range("D43").Value = "3.053.79"
$ws.Range("E43").Value = "  -8.89%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -8.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.276"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -10.56%  "
$ws.Range("D46").Value = "0.0₃0619"
$ws.Range("E46").Value = "  -17.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -11.43%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.25"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.90%  "
$ws.Range("E49").Value = "  -7.66%  "
$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.00"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.01%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.04"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.34%  "
